# Add a new "ignore" note column (E) next to the existing "Rien pour le
# moment..." placeholder cells in column D.
#
# For every row that currently has an empty D cell (style-only, no value)
# next to a C cell containing "Rien pour le moment...", we now:
#   - fill D with the same text/style as C ("Rien pour le moment...")
#   - add a new E cell, same style, containing the text "ignore"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(11, 12, 17, 24, 25, 26)

foreach ($r in $rows) {
    $src = $ws.Range("C$r")

    $dCell = $ws.Range("D$r")
    $src.Copy()
    $dCell.PasteSpecial(-4122)
    $dCell.Value = $src.Value2

    $eCell = $ws.Range("E$r")
    $src.Copy()
    $eCell.PasteSpecial(-4122)
    $eCell.Value = "ignore"
}

$ws.Range("E11:E12").Select()
